$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on price cells whose new value would otherwise be auto-detected as a number
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "62.937.34"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.544.31"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("D5").Value = "566.63"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "146.28"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").Value = "2.542.06"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "27.19"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").Value = "2.997.97"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "62.884.31"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "2.543.97"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "11.44"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "334.73"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "6.74"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("E26").Value = "  +5.14%  "
$ws.Range("E27").Value = "  +11.74%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "8.33"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  +6.47%  "
$ws.Range("D31").Value = "0.0₃0807"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "176.39"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("D35").Value = "407.59"
$ws.Range("E35").Value = "  +9.42%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "19.00"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "0.396"
$ws.Range("E37").Value = "  -1.03%  "
$ws.Range("D39").Value = "4.35"
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "39.05"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").Value = "152.70"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "20.72"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "0.604"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "0.0956"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").Value = "0.0517"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("D50").Value = "18.23"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  -0.87%  "
